# FaxAddressBookData.xlsx — 6thmay2020 commit
# Adds two new "query" worksheets (RecipientQuery, AddressBookQuery) that
# each hold a label cell ("Query") and a wrapped, multi-line SQL statement
# cell, and makes the last of the two sheets the active / selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: RecipientQuery
# ---------------------------------------------------------------------
$recipientSql = @"
SELECT 
    [FirstName] as 'Name1'
    ,[LastName] as 'Name2'
    ,[FaxNumber] as 'Fax Number'
    ,[LastChangedBy] as 'Last Changed By'
    ,[LastChangedOn] as 'Last Changed On'
    FROM [Fax_Recipient]
"@

$wsRecipient = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsRecipient.Name = "RecipientQuery"
$wsRecipient.Range("A1").Value = "Query"
$wsRecipient.Range("A2").Value = $recipientSql
$wsRecipient.Range("A2").WrapText = $true
$wsRecipient.Columns.Item(1).ColumnWidth = 55.166666666666664
$wsRecipient.Rows.Item(2).RowHeight = 105
[void]$wsRecipient.Range("A18").Select()

# ---------------------------------------------------------------------
# Sheet 2: AddressBookQuery
# ---------------------------------------------------------------------
$addressBookSql = @"
SELECT 
  [Name] as Name
  ,[FaxLine] as 'Fax Line'
  ,[LastChangedBy] as 'Last Changed By'
  ,[LastChangedOn] as 'Last Changed On'
  FROM [Fax_AddressBook]
"@

$wsAddressBook = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsAddressBook.Name = "AddressBookQuery"
$wsAddressBook.Range("A1").Value = "Query"
$wsAddressBook.Range("A2").Value = $addressBookSql
$wsAddressBook.Range("A2").WrapText = $true
$wsAddressBook.Columns.Item(1).ColumnWidth = 37.833333333333336
$wsAddressBook.Rows.Item(2).RowHeight = 90
[void]$wsAddressBook.Range("A10").Select()

# AddressBookQuery is the last sheet added, so it's already the active
# sheet/tab (mirrors the source workbook's activeTab="7"); re-activate
# explicitly for clarity / robustness.
$wsAddressBook.Activate()
